$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.711.44'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.641.95'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '216.35'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").Value = '0.0841'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").Value = '1.866.31'
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("D13").Value = '1.623.98'
$ws.Range("E13").Value = '  -2.98%  '
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("D17").Value = '26.714.33'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").Value = '214.15'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '2.42'
$ws.Range("E22").Value = '  +12.52%  '
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("E24").Value = '  -2.49%  '
$ws.Range("D25").Value = '145.66'
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").Value = '0.0509'
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("D33").Value = '3.00'
$ws.Range("E33").Value = '  -1.52%  '
$ws.Range("D34").Value = '1.292.15'
$ws.Range("E34").Value = '  +1.77%  '
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("D38").Value = '0.534'
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").Value = '0.818'
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").Value = '0.805'
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("E43").Value = '  -2.76%  '
$ws.Range("D44").Value = '1.791.24'
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("D45").Value = '61.37'
$ws.Range("E45").Value = '  +3.11%  '
$ws.Range("D46").Value = '91.24'
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '0.0526'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("E49").Value = '  -2.75%  '
$ws.Range("D50").Value = '7.66'
$ws.Range("E50").Value = '  -1.64%  '
$ws.Range("E51").Value = '  -0.16%  '
